$wb = $excel.ActiveWorkbook

# ALC row 5
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 620.06665
$ws.Range("I5").Value = 105.833336
$ws.Range("J5").Value = 962.8889
$ws.Range("K5").Value = 105.833336
$ws.Range("L5").Value = 962.8889
$ws.Range("M5").Value = 9.166663999999997
$ws.Range("N5").Value = -1192.8889

# ALC row 62
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 10224.5
$ws.Range("I62").Value = 10224.5
$ws.Range("K62").Value = 10224.5
$ws.Range("M62").Value = -9600.5

# ALC row 65
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H65").Value = 10224.5
$ws.Range("I65").Value = 10224.5
$ws.Range("K65").Value = 51122.5
$ws.Range("M65").Value = -48002.5

# ALC row 87
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H87").Value = 83328
$ws.Range("J87").Value = 83328
$ws.Range("L87").Value = 83328
$ws.Range("N87").Value = -85824

# ALC row 90
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H90").Value = 83328
$ws.Range("J90").Value = 83328
$ws.Range("L90").Value = 249984
$ws.Range("N90").Value = -262464

# ALC row 98
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 1589.3125
$ws.Range("I98").Value = 1589.3125
$ws.Range("J98").Value = 0
$ws.Range("K98").Value = 1589.3125
$ws.Range("L98").Value = 0
$ws.Range("M98").Value = -91.3125
$ws.Range("N98").ClearContents()

# ALC row 116
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H116").Value = 29281.572
$ws.Range("I116").Value = 0
$ws.Range("J116").Value = 29281.572
$ws.Range("K116").Value = 0
$ws.Range("L116").Value = 29281.572
$ws.Range("N116").Value = -36165.572
$ws.Range("M116").ClearContents()

# ALC row 122
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H122").Value = 1589.3125
$ws.Range("I122").Value = 1589.3125
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 4767.9375
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -2317.9375
$ws.Range("N122").ClearContents()

# ALC row 138
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 5337.1816
$ws.Range("J138").Value = 7941
$ws.Range("L138").Value = 23823
$ws.Range("N138").Value = -34103

# ALC row 141
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H141").Value = 2699.476
$ws.Range("J141").Value = 2992
$ws.Range("L141").Value = 8976
$ws.Range("N141").Value = -19336

# ARM row 15
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H15").Value = 6799.6665
$ws.Range("J15").Value = 9200
$ws.Range("L15").Value = 9200
$ws.Range("N15").Value = -9900

# ARM row 32
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4148.8354
$ws.Range("I32").Value = 4195.4873
$ws.Range("K32").Value = 4195.4873
$ws.Range("M32").Value = -3908.4873

# ARM row 45
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 8650.8125
$ws.Range("I45").Value = 8970.866
$ws.Range("K45").Value = 8970.866
$ws.Range("M45").Value = -8593.866

# ARM row 109
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H109").Value = 49699
$ws.Range("J109").Value = 49699
$ws.Range("L109").Value = 49699
$ws.Range("N109").Value = -52473

# ARM row 112
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H112").Value = 134399.8
$ws.Range("J112").Value = 134399.8
$ws.Range("L112").Value = 134399.8
$ws.Range("N112").Value = -137353.8

# ARM row 122
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 1670.2307
$ws.Range("I122").Value = 1577.6875
$ws.Range("J122").Value = 1818.3
$ws.Range("K122").Value = 4733.0625
$ws.Range("L122").Value = 5454.9
$ws.Range("M122").Value = -2283.0625
$ws.Range("N122").Value = -10354.9

# BSM row 22
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 517.25
$ws.Range("I22").Value = 524.8333
$ws.Range("J22").Value = 494.5
$ws.Range("K22").Value = 524.8333
$ws.Range("L22").Value = 494.5
$ws.Range("M22").Value = -351.8333
$ws.Range("N22").Value = -840.5

# CRP row 4
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 123049340
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 123049340
$ws.Range("K4").Value = 0
$ws.Range("L4").Value = 123049340
$ws.Range("N4").Value = -123049564
$ws.Range("M4").ClearContents()

# CRP row 8
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H8").Value = 9999
$ws.Range("I8").Value = 0
$ws.Range("J8").Value = 9999
$ws.Range("K8").Value = 0
$ws.Range("L8").Value = 9999
$ws.Range("N8").Value = -10279
$ws.Range("M8").ClearContents()

# CRP row 19
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H19").Value = 1637.75
$ws.Range("I19").Value = 765.4
$ws.Range("J19").Value = 5999.5
$ws.Range("K19").Value = 765.4
$ws.Range("L19").Value = 5999.5
$ws.Range("M19").Value = -595.4
$ws.Range("N19").Value = -6339.5

# CRP row 24
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H24").Value = 1637.75
$ws.Range("I24").Value = 765.4
$ws.Range("J24").Value = 5999.5
$ws.Range("K24").Value = 765.4
$ws.Range("L24").Value = 5999.5
$ws.Range("M24").Value = -595.4
$ws.Range("N24").Value = -6339.5

# CRP row 31
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2838.3262
$ws.Range("I31").Value = 1591.2593
$ws.Range("J31").Value = 4610.4736
$ws.Range("K31").Value = 1591.2593
$ws.Range("L31").Value = 4610.4736
$ws.Range("M31").Value = -1296.2593
$ws.Range("N31").Value = -5200.4736

# CRP row 34
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 2838.3262
$ws.Range("I34").Value = 1591.2593
$ws.Range("J34").Value = 4610.4736
$ws.Range("K34").Value = 1591.2593
$ws.Range("L34").Value = 4610.4736
$ws.Range("M34").Value = -1389.2593
$ws.Range("N34").Value = -5014.4736

# CRP row 99
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 8150.9165
$ws.Range("J99").Value = 10745.868
$ws.Range("L99").Value = 10745.868
$ws.Range("N99").Value = -13741.868

# CRP row 105
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value = 974
$ws.Range("I105").Value = 825.96
$ws.Range("K105").Value = 825.96
$ws.Range("M105").Value = 921.04

# CRP row 126
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H126").Value = 8150.9165
$ws.Range("J126").Value = 10745.868
$ws.Range("L126").Value = 32237.604
$ws.Range("N126").Value = -37177.604

# CRP row 141
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H141").Value = 251624.25
$ws.Range("I141").Value = 100000
$ws.Range("J141").Value = 259604.47
$ws.Range("K141").Value = 100000
$ws.Range("L141").Value = 259604.47
$ws.Range("M141").Value = -94820
$ws.Range("N141").Value = -269964.47

# CUL row 5
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 2083.3928
$ws.Range("I5").Value = 1209.3
$ws.Range("K5").Value = 3627.9
$ws.Range("M5").Value = -3515.9

# CUL row 8
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H8").Value = 709.5454999999999
$ws.Range("I8").Value = 709.5454999999999
$ws.Range("K8").Value = 2128.6365
$ws.Range("M8").Value = -1989.6365

# CUL row 11
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H11").Value = 139552.56
$ws.Range("I11").Value = 320.15585
$ws.Range("K11").Value = 960.46755
$ws.Range("M11").Value = -820.46755

# CUL row 129
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H129").Value = 2273.125
$ws.Range("J129").Value = 2780.4375
$ws.Range("L129").Value = 8341.3125
$ws.Range("N129").Value = -18341.3125

# CUL row 135
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H135").Value = 2083.3928
$ws.Range("I135").Value = 1209.3
$ws.Range("K135").Value = 10883.7
$ws.Range("M135").Value = -8348.699999999999

# LTW row 25
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H25").Value = 250000080
$ws.Range("I25").Value = 250000080
$ws.Range("K25").Value = 250000080
$ws.Range("M25").Value = -249999850

# LTW row 40
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 5111.909
$ws.Range("I40").Value = 5111.909
$ws.Range("K40").Value = 5111.909
$ws.Range("M40").Value = -4975.909

# LTW row 61
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 1556.1
$ws.Range("J61").Value = 1929.3334
$ws.Range("L61").Value = 1929.3334
$ws.Range("N61").Value = -2333.3334

# LTW row 82
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 2208.5908
$ws.Range("I82").Value = 2232.1667
$ws.Range("J82").Value = 2199.75
$ws.Range("K82").Value = 2232.1667
$ws.Range("L82").Value = 2199.75
$ws.Range("M82").Value = -1871.1667
$ws.Range("N82").Value = -2921.75

# LTW row 85
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H85").Value = 2208.5908
$ws.Range("I85").Value = 2232.1667
$ws.Range("J85").Value = 2199.75
$ws.Range("K85").Value = 2232.1667
$ws.Range("L85").Value = 2199.75
$ws.Range("M85").Value = -984.1667000000002
$ws.Range("N85").Value = -4695.75

# LTW row 113
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H113").Value = 1556.1
$ws.Range("J113").Value = 1929.3334
$ws.Range("L113").Value = 1929.3334
$ws.Range("N113").Value = -6269.3334

# LTW row 122
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 5095.952
$ws.Range("I122").Value = 4317.5
$ws.Range("K122").Value = 12952.5
$ws.Range("M122").Value = -10502.5

# WVR row 24
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H24").Value = 19800
$ws.Range("J24").Value = 0
$ws.Range("L24").Value = 0
$ws.Range("N24").ClearContents()

# WVR row 122
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 3720.4
$ws.Range("I122").Value = 3866
$ws.Range("J122").Value = 3593
$ws.Range("K122").Value = 11598
$ws.Range("L122").Value = 10779
$ws.Range("M122").Value = -9148
$ws.Range("N122").Value = -15679
